$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates of Price (D) and Volume(1h) (E) columns,
# plus the three-row reorder of FraxShare / ARBITRUM / HuobiToken (rows 46-48).
# Cells whose new text would otherwise be auto-coerced to a Number by Excel
# (e.g. "22.30" -> 22.3, dropping the trailing zero) are forced to stay text
# by temporarily applying a "@" (Text) number format, then restoring the
# original style so no stray formatting is left behind.

$ws.Range('D2').Value = '39.489.77'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').Value = '2.163.90'
$ws.Range('E3').Value = '  +3.62%  '
$ws.Range('E4').Value = '  +0.02%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.92'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  +1.36%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.07'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +6.82%  '
$ws.Range('E8').Value = '  +0.01%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.399'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +3.92%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0864'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +2.95%  '
$ws.Range('E11').Value = '  -0.15%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.96'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +6.50%  '
$ws.Range('D13').Value = '2.484.13'
$ws.Range('E13').Value = '  +3.54%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.30'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +2.25%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.818'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +2.42%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.57'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').Value = '2.169.51'
$ws.Range('E17').Value = '  +3.91%  '
$ws.Range('D18').Value = '39.447.02'
$ws.Range('E18').Value = '  +2.10%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.36'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +1.26%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.16'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('D21').Value = '0.0₃0853'
$ws.Range('E21').Value = '  +1.91%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '232.74'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('E23').Value = '  -0.01%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +1.99%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -2.82%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.64'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +1.11%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.52'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +1.00%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.138'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.38%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.91'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +3.97%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.41'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -3.95%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.66'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +11.71%  '
$ws.Range('E32').Value = '  +1.70%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.68'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +4.09%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.85'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +3.95%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.10'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +10.25%  '
$ws.Range('E36').Value = '  +2.67%  '
$ws.Range('E37').Value = '  +1.67%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.65'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +3.20%  '
$ws.Range('E39').Value = '  +0.01%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '105.18'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +5.15%  '
$ws.Range('E41').Value = '  +1.14%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.93'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').Value = '1.540.53'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('E44').Value = '  +6.69%  '
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.11'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +7.81%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.82'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.83'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('E49').Value = '  +3.60%  '
$ws.Range('D50').Value = '2.369.58'
$ws.Range('E50').Value = '  +3.58%  '
$ws.Range('E51').Value = '  +0.22%  '
